$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SummaryReport")

$ws.Range("E2").Value = 32.417999999999999
$ws.Range("F2").Value = 6.4930000000000003
$ws.Range("G2").Value = 31.456

$ws.Range("E3").Value = 0.107
$ws.Range("F3").Value = 0.012
$ws.Range("G3").Value = 0.081000000000000003

$ws.Range("E4").Value = 0.13100000000000001
$ws.Range("F4").Value = 0.014
$ws.Range("G4").Value = 0.078

$ws.Range("E5").Value = 0.97899999999999998
$ws.Range("F5").Value = 0.189
$ws.Range("G5").Value = 0.72699999999999998

$ws.Range("E6").Value = 1.417
$ws.Range("F6").Value = 0.14499999999999999
$ws.Range("G6").Value = 0.83099999999999996

$ws.Range("E7").Value = 0.19500000000000001
$ws.Range("F7").Value = 0.017999999999999999
$ws.Range("G7").Value = 0.154

$ws.Range("E8").Value = 0.17499999999999999
$ws.Range("F8").Value = 0.017000000000000001
$ws.Range("G8").Value = 0.13100000000000001

$ws.Range("E9").Value = 0.23300000000000001
$ws.Range("F9").Value = 0.023
$ws.Range("G9").Value = 0.153

$ws.Range("E10").Value = 0.14299999999999999
$ws.Range("F10").Value = 0.014
$ws.Range("G10").Value = 0.106

$ws.Range("E11").Value = 0.078
$ws.Range("F11").Value = 0.0089999999999999993

$ws.Range("E12").Value = 0.24199999999999999
$ws.Range("F12").Value = 0.02
$ws.Range("G12").Value = 0.153

$ws.Range("E13").Value = 0.217
$ws.Range("F13").Value = 0.025999999999999999
$ws.Range("G13").Value = 0.159

$ws.Range("E14").Value = 0.158
$ws.Range("F14").Value = 0.014999999999999999
$ws.Range("G14").Value = 0.13700000000000001

$ws.Range("E15").Value = 10.448
$ws.Range("F15").Value = 0.028000000000000001
$ws.Range("G15").Value = 10.428000000000001

$ws.Range("E16").Value = 25.908999999999999
$ws.Range("F16").Value = 0.049000000000000002
$ws.Range("G16").Value = 25.78

$ws.Range("E17").Value = 32.417000000000002
$ws.Range("F17").Value = 0.21199999999999999
$ws.Range("G17").Value = 31.555

$ws.Range("E18").Value = 21.376999999999999
$ws.Range("F18").Value = 0.062
$ws.Range("G18").Value = 21.366

$ws.Range("E19").Value = 22.47
$ws.Range("F19").Value = 0.26300000000000001
$ws.Range("G19").Value = 21.87

$ws.Range("E20").Value = 16.777000000000001
$ws.Range("F20").Value = 0.055
$ws.Range("G20").Value = 16.663

$ws1 = $wb.Worksheets.Item("Автоматизированный расчет")
$ws1.Range("C54").Select()
